# Updates cryptos list: refresh Price and Volume(1h) columns, and
# shift Coin/Link pairs to reflect the latest coinranking.com ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.465.66'
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("D3").Value = '2.351.90'
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '''521.98'
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("D6").Value = '''137.16'
$ws.Range("E6").Value = '  +3.31%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '''0.539'
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '''0.103'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").Value = '''5.44'
$ws.Range("E10").Value = '  +5.77%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '''0.152'
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '''0.344'
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = '''24.15'
$ws.Range("E13").Value = '  +0.69%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.772.02'
$ws.Range("E14").Value = '  +1.72%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '57.438.05'
$ws.Range("E15").Value = '  +2.15%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000136'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.369.27'
$ws.Range("E17").Value = '  +2.20%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '''10.68'
$ws.Range("E18").Value = '  +1.55%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''329.12'
$ws.Range("E19").Value = '  +2.51%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '''4.26'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''6.75'
$ws.Range("E21").Value = '  +1.24%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''61.45'
$ws.Range("E23").Value = '  +1.34%  '
$ws.Range("B24").Value = 'Kaspa'
$ws.Range("C24").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D24").Value = '''0.166'
$ws.Range("E24").Value = '  +4.99%  '
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").Value = '''0.995'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '''8.24'
$ws.Range("E26").Value = '  +8.05%  '
$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D27").Value = '''1.32'
$ws.Range("E27").Value = '  +10.83%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0751'
$ws.Range("E28").Value = '  +3.07%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '''170.37'
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''1.71'
$ws.Range("E30").Value = '  -0.52%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = '''6.29'
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '''18.64'
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''1.30'
$ws.Range("E34").Value = '  +2.88%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = '''0.994'
$ws.Range("E35").Value = '  -0.36%  '
$ws.Range("B36").Value = 'SuiNetwork'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D36").Value = '''0.928'
$ws.Range("E36").Value = '  +0.35%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '''4.06'
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '''1.59'
$ws.Range("E38").Value = '  +5.14%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '''38.56'
$ws.Range("E39").Value = '  +3.11%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '''151.07'
$ws.Range("E40").Value = '  +8.60%  '
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").Value = '''0.386'
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '''3.66'
$ws.Range("E42").Value = '  +2.34%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''5.33'
$ws.Range("E43").Value = '  +4.97%  '
$ws.Range("D44").Value = '''283.92'
$ws.Range("E44").Value = '  +4.79%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '''0.0940'
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '''0.0509'
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '''0.566'
$ws.Range("E47").Value = '  +1.97%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '''0.0221'
$ws.Range("E48").Value = '  +2.69%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '''18.19'
$ws.Range("E49").Value = '  +6.94%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''17.74'
$ws.Range("E50").Value = '  +5.02%  '
$ws.Range("B51").Value = 'Polygon'
$ws.Range("C51").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D51").Value = '''0.384'
$ws.Range("E51").Value = '  +0.72%  '
